$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "dSF" column (F) values for the rows that were repulled/recalculated
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("F14").Value = 3
